$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $v = $cell.Value2
    if ($v -eq 45174) {
        $cell.Value = 45175
    }
}
